# regen sval data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(0.06328177979961902, 0.05231270169004087, 0.1529057820181812, 0.4998867070740569, 0.768386970581898)
    3 = @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 6.048734245549538)
    4 = @(0.3464964993005633, 0.3375848360084654, 3.082599426703578, 0.4998867070740569, 4.266567469086664)
    5 = @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 6.48142807727062, 12.0302756157461)
    6 = @(0.3464964993005633, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 2.652525447291612)
    7 = @(0.3464964993005633, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 2.652525447291612)
    8 = @(3.182878228561681, 1.65323645889881, 3.082599426703578, 0.4998867070740569, 8.418600821238126)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 7).Value = $vals[4]
}
